# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (fund holding detail) right before the
#    "总计" (totals) sheet, copying the header / index-column formatting
#    from the most recent quarter sheet ("2021-Q4").
# 2) Insert a new top data row into "总计" for "2022-Q1" (7 funds held,
#    3.32 billion yuan), shifting the older rows down and renumbering the
#    index column.
#
# NOTE: worksheet object references resolve by tab *position*, not a
# fixed identity, so every handle is re-fetched by name right before each
# use (instead of being reused across the sheet-insert below).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet just before "总计"
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Re-fetch by name: the position-based handles above are now stale since
# the insert shifted "总计" one slot to the right.
$newSheet = $wb.Worksheets.Item("2022-Q1")
$refSheet = $wb.Worksheets.Item("2021-Q4")

# Pull the header-row and index-column formatting from the reference sheet
# so the new sheet matches the look of the other quarterly sheets.
$refSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$refSheet.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding rows. Columns B-G are kept as text (fund codes with leading
# zeros, and numeric-looking percentages/amounts stored as text) by
# stamping a text number format before the write, then restoring the
# default "Normal" style so no residual format id lingers on the cell.
$rows = @(
  @("011866", "广发价值增长混合型证券投资基金A", "18.06", "92.01", "6.18", "1.1161", 7),
  @("002624", "广发优企精选灵活配置混合A", "12.98", "92.40", "6.19", "0.8035", 8),
  @("270025", "广发行业领先混合A", "11.11", "91.67", "5.75", "0.6388", 8),
  @("960001", "广发行业领先混合H", "11.11", "91.67", "5.75", "0.6388", 8),
  @("000747", "广发逆向策略灵活配置混合", "1.25", "89.61", "4.96", "0.0620", 9),
  @("011867", "广发价值增长混合型证券投资基金C", "0.81", "92.01", "6.18", "0.0501", 7),
  @("010021", "广发优企精选灵活配置混合C", "0.15", "92.40", "6.19", "0.0093", 8)
)

$r = 2
foreach ($row in $rows) {
  $newSheet.Cells.Item($r, 1).Value = ($r - 2)

  # Columns B..G (2..7) are fund code / name / size / position / ratio /
  # holding value - all stored as text in the source data.
  for ($c = 2; $c -le 7; $c++) {
    $cell = $newSheet.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $row[$c - 2]
    $cell.Style = "Normal"
  }

  # Column H (8) is the numeric position rank.
  $newSheet.Cells.Item($r, 8).Value = $row[6]

  $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: insert the "2022-Q1" summary row at the top of "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row picks up formatting from its neighbours; reset it to
# the plain/default style before laying down the real values.
$totalSheet.Range("A2:D2").Style = "Normal"

$refSheet = $wb.Worksheets.Item("2021-Q4")
$refSheet.Range("A2").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 3.32

# Renumber the index column for the rows that got pushed down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
